$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Em dash character used as placeholder text in several cells
$emDash = [char]0x2014

# ---------------------------------------------------------------
# Row 16 changes
# ---------------------------------------------------------------
# A16: was a datetime-like string "2025-11-07 00:00:00", now literal text "11/7/2025"
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "11/7/2025"

# Q16: "No" -> empty
$ws.Range("Q16").ClearContents()

# Y16: empty -> em dash
$ws.Range("Y16").Value = $emDash

# AF16: empty -> em dash
$ws.Range("AF16").Value = $emDash

# ---------------------------------------------------------------
# Row 18 changes
# ---------------------------------------------------------------
# A18: empty -> "Invalid Date"
$ws.Range("A18").Value = "Invalid Date"

# L18: "Immediate" -> empty
$ws.Range("L18").ClearContents()

# O18: empty -> em dash
$ws.Range("O18").Value = $emDash

# P18: empty -> em dash
$ws.Range("P18").Value = $emDash

# Q18: "No" -> empty
$ws.Range("Q18").ClearContents()

# W18: empty -> em dash
$ws.Range("W18").Value = $emDash

# Y18: empty -> em dash
$ws.Range("Y18").Value = $emDash

# AF18: empty -> em dash
$ws.Range("AF18").Value = $emDash
